$wb = $excel.ActiveWorkbook

# ---- Sheet ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 63856.438
$ws.Range("I70").Value = 201060.0
$ws.Range("J70").Value = 1491.1818
$ws.Range("K70").Value = 603180.0
$ws.Range("L70").Value = 4473.5454
$ws.Range("M70").Value = -602910.0
$ws.Range("N70").Value = -5013.5454
$ws.Range("H73").Value = 63856.438
$ws.Range("I73").Value = 201060.0
$ws.Range("J73").Value = 1491.1818
$ws.Range("K73").Value = 603180.0
$ws.Range("L73").Value = 4473.5454
$ws.Range("M73").Value = -602244.0
$ws.Range("N73").Value = -6345.5454
$ws.Range("H76").Value = 4999.3
$ws.Range("I76").Value = 5211.625
$ws.Range("J76").Value = 4150.0
$ws.Range("K76").Value = 5211.625
$ws.Range("L76").Value = 4150.0
$ws.Range("M76").Value = -4896.625
$ws.Range("N76").Value = -4780.0
$ws.Range("H79").Value = 4999.3
$ws.Range("I79").Value = 5211.625
$ws.Range("J79").Value = 4150.0
$ws.Range("K79").Value = 5211.625
$ws.Range("L79").Value = 4150.0
$ws.Range("M79").Value = -4119.625
$ws.Range("N79").Value = -6334.0
$ws.Range("H80").Value = 8724.429
$ws.Range("I80").Value = 937.3333
$ws.Range("J80").Value = 22741.2
$ws.Range("K80").Value = 2811.9999
$ws.Range("L80").Value = 68223.6
$ws.Range("M80").Value = -1813.9999
$ws.Range("N80").Value = -70219.6
$ws.Range("H83").Value = 8724.429
$ws.Range("I83").Value = 937.3333
$ws.Range("J83").Value = 22741.2
$ws.Range("K83").Value = 8435.9997
$ws.Range("L83").Value = 204670.8
$ws.Range("M83").Value = -3443.9997
$ws.Range("N83").Value = -214654.8
$ws.Range("H112").Value = 2528.9863
$ws.Range("I112").Value = 250.0
$ws.Range("J112").Value = 2661.1016
$ws.Range("K112").Value = 750.0
$ws.Range("L112").Value = 7983.3048
$ws.Range("M112").Value = 358.0
$ws.Range("N112").Value = -10199.3048
$ws.Range("H116").Value = 6252283.5
$ws.Range("I116").Value = 15386546.0
$ws.Range("J116").Value = 2524.7368
$ws.Range("K116").Value = 15386546.0
$ws.Range("L116").Value = 2524.7368
$ws.Range("M116").Value = -15383104.0
$ws.Range("N116").Value = -9408.7368
$ws.Range("H132").Value = 2042.4746
$ws.Range("I132").Value = 1435.5
$ws.Range("J132").Value = 3320.3157
$ws.Range("K132").Value = 4306.5
$ws.Range("L132").Value = 9960.947100000001
$ws.Range("M132").Value = -1776.5
$ws.Range("N132").Value = -15020.9471
$ws.Range("H138").Value = 2217.6309
$ws.Range("I138").Value = 1241.5862
$ws.Range("J138").Value = 4394.9614
$ws.Range("K138").Value = 3724.7586
$ws.Range("L138").Value = 13184.8842
$ws.Range("M138").Value = 1415.2414
$ws.Range("N138").Value = -23464.8842
$ws.Range("H139").Value = 41900.0
$ws.Range("J139").Value = 41900.0
$ws.Range("L139").Value = 41900.0
$ws.Range("N139").Value = -52180.0
$ws.Range("H141").Value = 3460.1897
$ws.Range("I141").Value = 1544.0577
$ws.Range("J141").Value = 20066.666
$ws.Range("K141").Value = 4632.1731
$ws.Range("L141").Value = 60199.99800000001
$ws.Range("M141").Value = 547.8269
$ws.Range("N141").Value = -70559.998

# ---- Sheet ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 7616.086
$ws.Range("I32").Value = 7414.906
$ws.Range("K32").Value = 7414.906
$ws.Range("M32").Value = -7127.906
$ws.Range("H49").Value = 6000.0
$ws.Range("I49").Value = 10000.0
$ws.Range("J49").Value = 2000.0
$ws.Range("K49").Value = 10000.0
$ws.Range("L49").Value = 2000.0
$ws.Range("M49").Value = -9740.0
$ws.Range("N49").Value = -2520.0
$ws.Range("H61").Value = 2116.575
$ws.Range("I61").Value = 1110.3715
$ws.Range("J61").Value = 9160.0
$ws.Range("K61").Value = 1110.3715
$ws.Range("L61").Value = 9160.0
$ws.Range("M61").Value = -898.3715
$ws.Range("N61").Value = -9584.0
$ws.Range("H97").Value = 1184.909
$ws.Range("I97").Value = 961.25
$ws.Range("K97").Value = 961.25
$ws.Range("M97").Value = -465.25
$ws.Range("H136").Value = 2116.575
$ws.Range("I136").Value = 1110.3715
$ws.Range("J136").Value = 9160.0
$ws.Range("K136").Value = 3331.1145
$ws.Range("L136").Value = 27480.0
$ws.Range("M136").Value = -781.1144999999997
$ws.Range("N136").Value = -32580.0
$ws.Range("H139").Value = 58338.934
$ws.Range("J139").Value = 58006.0
$ws.Range("L139").Value = 58006.0
$ws.Range("N139").Value = -68286.0

# ---- Sheet BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2631.1082
$ws.Range("I134").Value = 2279.963
$ws.Range("K134").Value = 6839.889000000001
$ws.Range("M134").Value = -4304.889000000001

# ---- Sheet CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H3").Value = 35000.0
$ws.Range("I3").Value = 50000.0
$ws.Range("J3").Value = 20000.0
$ws.Range("K3").Value = 50000.0
$ws.Range("L3").Value = 20000.0
$ws.Range("M3").Value = -49887.0
$ws.Range("N3").Value = -20226.0
$ws.Range("H7").Value = 219.9
$ws.Range("I7").Value = 166.0
$ws.Range("J7").Value = 273.8
$ws.Range("K7").Value = 166.0
$ws.Range("L7").Value = 273.8
$ws.Range("M7").Value = -53.0
$ws.Range("N7").Value = -499.8
$ws.Range("H31").Value = 1624.8246
$ws.Range("I31").Value = 1268.4166
$ws.Range("K31").Value = 1268.4166
$ws.Range("M31").Value = -973.4166
$ws.Range("H34").Value = 1624.8246
$ws.Range("I34").Value = 1268.4166
$ws.Range("K34").Value = 1268.4166
$ws.Range("M34").Value = -1066.4166
$ws.Range("H122").Value = 2992.2083
$ws.Range("I122").Value = 2985.6924
$ws.Range("J122").Value = 2999.9092
$ws.Range("K122").Value = 8957.0772
$ws.Range("L122").Value = 8999.7276
$ws.Range("M122").Value = -6507.0772
$ws.Range("N122").Value = -13899.7276
$ws.Range("H132").Value = 399253.84
$ws.Range("I132").Value = 484021.72
$ws.Range("K132").Value = 1452065.16
$ws.Range("M132").Value = -1449535.16

# ---- Sheet CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 24.083334
$ws.Range("I14").Value = 24.083334
$ws.Range("K14").Value = 72.250002
$ws.Range("M14").Value = 100.749998
$ws.Range("H69").Value = 2000.0
$ws.Range("I69").Value = 0.0
$ws.Range("J69").Value = 2000.0
$ws.Range("K69").Value = 0.0
$ws.Range("L69").Value = 6000.0
$ws.Range("M69").ClearContents()
$ws.Range("N69").Value = -7622.0
$ws.Range("H72").Value = 2000.0
$ws.Range("I72").Value = 0.0
$ws.Range("J72").Value = 2000.0
$ws.Range("K72").Value = 0.0
$ws.Range("L72").Value = 18000.0
$ws.Range("M72").ClearContents()
$ws.Range("N72").Value = -26112.0
$ws.Range("H131").Value = 998.03
$ws.Range("J131").Value = 1021.59375
$ws.Range("L131").Value = 3064.78125
$ws.Range("N131").Value = -13144.78125
$ws.Range("H133").Value = 5434.615
$ws.Range("I133").Value = 2732.75
$ws.Range("J133").Value = 6635.4443
$ws.Range("K133").Value = 8198.25
$ws.Range("L133").Value = 19906.3329
$ws.Range("M133").Value = -3138.25
$ws.Range("N133").Value = -30026.3329
$ws.Range("H134").Value = 4868.12
$ws.Range("I134").Value = 3173.0
$ws.Range("J134").Value = 5998.2
$ws.Range("K134").Value = 9519.0
$ws.Range("L134").Value = 17994.6
$ws.Range("M134").Value = -4449.0
$ws.Range("N134").Value = -28134.6
$ws.Range("H138").Value = 2479.1538
$ws.Range("I138").Value = 914.8333
$ws.Range("K138").Value = 2744.4999
$ws.Range("M138").Value = 2395.5001

# ---- Sheet GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H43").Value = 4779.25
$ws.Range("I43").Value = 2017.0
$ws.Range("J43").Value = 5700.0
$ws.Range("K43").Value = 2017.0
$ws.Range("L43").Value = 5700.0
$ws.Range("M43").Value = -1866.0
$ws.Range("N43").Value = -6002.0
$ws.Range("H62").Value = 30000.0
$ws.Range("J62").Value = 30000.0
$ws.Range("L62").Value = 30000.0
$ws.Range("N62").Value = -31372.0
$ws.Range("H65").Value = 30000.0
$ws.Range("J65").Value = 30000.0
$ws.Range("L65").Value = 90000.0
$ws.Range("N65").Value = -96864.0
$ws.Range("H80").Value = 2848.9333
$ws.Range("I80").Value = 2894.875
$ws.Range("J80").Value = 2665.1667
$ws.Range("K80").Value = 2894.875
$ws.Range("L80").Value = 2665.1667
$ws.Range("M80").Value = -1896.875
$ws.Range("N80").Value = -4661.1667
$ws.Range("H83").Value = 2848.9333
$ws.Range("I83").Value = 2894.875
$ws.Range("J83").Value = 2665.1667
$ws.Range("K83").Value = 14474.375
$ws.Range("L83").Value = 13325.8335
$ws.Range("M83").Value = -9482.375
$ws.Range("N83").Value = -23309.8335

# ---- Sheet LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1601.2
$ws.Range("J16").Value = 2000.6666
$ws.Range("L16").Value = 2000.6666
$ws.Range("N16").Value = -2340.6666
$ws.Range("H34").Value = 50000.0
$ws.Range("I34").Value = 50000.0
$ws.Range("J34").Value = 0.0
$ws.Range("K34").Value = 50000.0
$ws.Range("L34").Value = 0.0
$ws.Range("M34").Value = -49828.0
$ws.Range("N34").ClearContents()
$ws.Range("H56").Value = 23541.25
$ws.Range("I56").Value = 15210.2
$ws.Range("J56").Value = 37426.332
$ws.Range("K56").Value = 15210.2
$ws.Range("L56").Value = 37426.332
$ws.Range("M56").Value = -14519.2
$ws.Range("N56").Value = -38808.332
$ws.Range("H132").Value = 5147.1333
$ws.Range("I132").Value = 4954.8184
$ws.Range("K132").Value = 14864.4552
$ws.Range("M132").Value = -12334.4552
$ws.Range("H136").Value = 2209.25
$ws.Range("I136").Value = 1762.5135
$ws.Range("J136").Value = 4570.5713
$ws.Range("K136").Value = 5287.5405
$ws.Range("L136").Value = 13711.7139
$ws.Range("M136").Value = -2737.5405
$ws.Range("N136").Value = -18811.7139

# ---- Sheet WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H14").Value = 27502.0
$ws.Range("I14").Value = 27502.0
$ws.Range("J14").Value = 0.0
$ws.Range("K14").Value = 27502.0
$ws.Range("L14").Value = 0.0
$ws.Range("M14").Value = -27334.0
$ws.Range("N14").ClearContents()
$ws.Range("H107").Value = 629.619
$ws.Range("I107").Value = 624.64703
$ws.Range("J107").Value = 650.75
$ws.Range("K107").Value = 1873.94109
$ws.Range("L107").Value = 1952.25
$ws.Range("M107").Value = 46.0589100000002
$ws.Range("N107").Value = -5792.25
$ws.Range("H132").Value = 1326.0
$ws.Range("I132").Value = 900.3333
$ws.Range("K132").Value = 2700.9999
$ws.Range("M132").Value = -170.9998999999998
$ws.Range("H136").Value = 1002.64105
$ws.Range("I136").Value = 1046.4324
$ws.Range("J136").Value = 192.5
$ws.Range("K136").Value = 3139.2972
$ws.Range("L136").Value = 577.5
$ws.Range("M136").Value = -589.2972
$ws.Range("N136").Value = -5677.5
